$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing current rows 4-6 down to 5-7.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the data for the inserted record.
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(4, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(4, 4).Value = 44544
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
$ws.Cells.Item(4, 5).Value = 15
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100103
$ws.Cells.Item(4, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(4, 9).Value = 100103003
$ws.Cells.Item(4, 10).Value = "Damasco"
$ws.Cells.Item(4, 11).Value = "Castle Brite"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 250
$ws.Cells.Item(4, 14).Value = 20000
$ws.Cells.Item(4, 15).Value = 22000
$ws.Cells.Item(4, 16).Value = 21000
$ws.Cells.Item(4, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(4, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(4, 19).Value = 1167
$ws.Cells.Item(4, 20).Value = 18
